$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Cells.Item(272, 1).Value = "14 17:26>>> 4FD889D140   Freddy Velez"
$ws.Cells.Item(273, 1).Value = "15 14:46>>> CFD89370C0   John Tomanelli"
$ws.Cells.Item(274, 1).Value = "15 14:47>>> CFD89370C0   John Tomanelli"
